$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at position 10 (shifts old rows 10-11 down to 11-12)
$ws.Rows("10:10").Insert()
$ws.Rows("10:10").RowHeight = 42.75

# 2) Populate the new row's cells. The order below controls the order new
#    strings are appended to the shared string table (matches the source file).
$ws.Range("C10").Value = "J4"
$ws.Range("E10").Value = "Connector Header Through Hole, Right Angle 4 position 0.165"" (4.20mm)"
$ws.Range("G10").Value = "Connector_Molex:Molex_Mini-Fit_Jr_5569-04A2_2x02_P4.20mm_Horizontal"
$ws.Range("F10").Value = "https://www.molex.com/en-us/products/part-detail/26013115?display=pdf"
$ws.Range("H10").Value = "https://www.digikey.jp/en/products/detail/molex/0026013115/4119889"
$ws.Range("D10").Value = "Molex 5569"
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 1
$ws.Range("I10").Value = 26013115

# 3) Fix up the style of C10 (Reference column) to match the other
#    "centered bold-ish" columns on this row (copy format from D9).
$ws.Range("D9").Copy()
$ws.Range("C10").PasteSpecial(-4122)

# 4) Rebuild all hyperlinks in the correct final order/targets. This engine
#    keeps hyperlink anchors fixed to their original cell refs when rows are
#    inserted/shifted, so we clear everything and re-add them so that the
#    four links belonging to the shifted rows point at their new cells, and
#    two new links are appended for the newly inserted row.
$ws.Range("H10").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("H8"), "https://www.digikey.jp/en/products/detail/amphenol-cs-commercial-products/RJHSE-5380/1242692")
$ws.Hyperlinks.Add($ws.Range("H5"), "https://www.digikey.jp/en/products/detail/rubycon/50ZLH100MEFC8X11-5/3563386")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.assmann-wsw.com/uploads/datasheets/ASS_4888_CO.pdf")
$ws.Hyperlinks.Add($ws.Range("H2"), "https://www.digikey.jp/en/products/detail/samsung-electro-mechanics/CL10A105KA8NNNC/3886760")
$ws.Hyperlinks.Add($ws.Range("H3"), "https://www.digikey.jp/en/products/detail/samsung-electro-mechanics/CL10B104KB8NNWC/3887593")
$ws.Hyperlinks.Add($ws.Range("H4"), "https://www.digikey.jp/en/products/detail/kemet/C0603C102K4RECAUTO/8640472")
$ws.Hyperlinks.Add($ws.Range("H9"), "https://www.digikey.jp/en/products/detail/assmann-wsw-components/A-DS-09-A-KG-T2S/1241804")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://mm.digikey.com/Volume0/opasdata/d220001/medias/docus/609/CL10A105KA8NNNC_Spec.pdf")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://mm.digikey.com/Volume0/opasdata/d220001/medias/docus/658/CL10B104KB8NNWC_Spec.pdf")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.rubycon.co.jp/wp-content/uploads/catalog-aluminum/ZLH.pdf")
$ws.Hyperlinks.Add($ws.Range("H6"), "https://www.digikey.jp/en/products/detail/liteon/LTST-C191KGKT/386835")
$ws.Hyperlinks.Add($ws.Range("H7"), "https://akizukidenshi.com/catalog/g/g106282/")
$ws.Hyperlinks.Add($ws.Range("H11"), "https://www.digikey.jp/en/products/detail/torex-semiconductor-ltd/XC6216D332PR-G/2815608?s=N4IgTCBcDaIBoGEBsYCMSAmBmLYAOATgLQDmIAugL5A")
$ws.Hyperlinks.Add($ws.Range("H12"), "https://www.digikey.jp/en/products/detail/texas-instruments/SN65LVDS2DR/1574870")
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.ti.com/lit/ds/symlink/sn65lvds2.pdf")
$ws.Hyperlinks.Add($ws.Range("F11"), "https://product.torexsemi.com/system/files/series/xc6216.pdf")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://cdn.amphenol-cs.com/media/wysiwyg/files/drawing/rjhsex380.pdf")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.semitec-global.com/uploads/2022/01/P22-23-CRD.pdf")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://optoelectronics.liteon.com/upload/download/DS22-2000-228/LTST-C191KGKT.PDF")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://connect.kemet.com:7667/gateway/IntelliData-ComponentDocumentation/1.0/download/datasheet/C0603C102K4RECAUTO")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.molex.com/en-us/products/part-detail/26013115?display=pdf")
$ws.Hyperlinks.Add($ws.Range("H10"), "https://www.digikey.jp/en/products/detail/molex/0026013115/4119889")

# 5) Update the sheet view (scroll position / active selection)
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("I15").Select()
